# Stundenliste update: add new entries (Besprechung / Besprechung, erste
# Design Schritte / template_app work, then later Design Class Diagram /
# Refinement Class Diagram) around the existing "Recherche..." row and the
# "Decoding first UART Signal" row, and extend the running-total formulas
# to match the new, shorter tail of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 3 fresh rows right above the old row 26 ("Decoding first UART
# Signal"), which pushes it (and everything below) down to row 29 and
# carries that row's formatting (date format, row height, wrap text) along
# with it into the new rows. ------------------------------------------------
$ws.Rows("26:28").Insert()

# The sheet used to run to row 36; the new layout only needs to run to row
# 33, so drop the now-surplus trailing rows (the old rows 34:36 plus the
# three blank carry-over rows created by the insert above).
$ws.Rows("34:39").Delete()

# --- Fill in the new rows with the actual entries --------------------------
$ws.Range("A26").Value = (Get-Date -Year 2021 -Month 1 -Day 19).Date
$ws.Range("B26").Value = 1
$ws.Range("D26").Value = "Besprechung"

$ws.Range("A27").Value = (Get-Date -Year 2021 -Month 1 -Day 21).Date
$ws.Range("B27").Value = 5
$ws.Range("D27").Value = "Besprechung, erste Design Schritte"

$ws.Range("A28").Value = (Get-Date -Year 2021 -Month 1 -Day 23).Date
$ws.Range("B28").Value = 3
$ws.Range("D28").Value = "template_app, Makfiles"

# Row 29 already holds the shifted-down "Decoding first UART Signal" entry
# (date 44220, 9 hours) - nothing to change there.

# Rows 30/31 were previously blank (General format) in column A, so copy the
# date cell formatting from row 25 first (keeps the same date-format style
# index instead of Excel minting a brand-new numFmt for them).
$ws.Range("A25").Copy($ws.Range("A30"))
$ws.Range("A30").Value = 44221
$ws.Range("B30").Value = 5
$ws.Range("D30").Value = "Design Class Diagram"

$ws.Range("A25").Copy($ws.Range("A31"))
$ws.Range("A31").Value = 44223
$ws.Range("B31").Value = 1.5
$ws.Range("D31").Value = "Refinement Class Diagram"

# Rows 32 and 33 stay blank in columns A/B/D, same as before.

# --- Re-establish the running-total formula down through row 33 -----------
for ($r = 26; $r -le 33; $r++) {
    $prev = $r - 1
    $ws.Range("C$r").Formula = "=C$prev+B$r"
}

# Leave the cursor where the author last left it.
$ws.Range("C31").Select()
